# Update the AppControl sheet: replace the long email list in B25 with a
# single address, and turn that cell into a mailto: hyperlink.
$wb = $excel.ActiveWorkbook

$wsAppControl = $wb.Worksheets.Item("AppControl")
$wsAppControl.Range("B25").Value = "stiyyagura@enhops.com"
$wsAppControl.Hyperlinks.Add($wsAppControl.Range("B25"), "mailto:stiyyagura@enhops.com")

$wsAppControl.Range("B27").Select()

# Update the smoke sheet: flip the run flag for rows 3-18 from Y to N, and
# move the active selection/scroll position.
$wsSmoke = $wb.Worksheets.Item("smoke")
for ($r = 3; $r -le 18; $r++) {
    $wsSmoke.Range("B$r").Value = "N"
}

$wsSmoke.Range("B21").Select()
